$d = $word.ActiveDocument

# --- Portuguese "Programa" paragraph (detailed) ---
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Campos de atuação. A Física como", $true, $false, $false, $false, $false, $true, 1, $false, "Campos de atuação. ^lA Física como", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("significativos de Física.Conceitos básicos", $true, $false, $false, $false, $false, $true, 1, $false, "significativos de Física.^lConceitos básicos", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("de um engenheiro.Desenvolvimento de um projeto temático", $true, $false, $false, $false, $false, $true, 1, $false, "de um engenheiro.^lDesenvolvimento de um projeto temático", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Engenharia Física.Competição entre projetos", $true, $false, $false, $false, $false, $true, 1, $false, "Engenharia Física.^lCompetição entre projetos", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("diferentes grupos.Avaliação das competições", $true, $false, $false, $false, $false, $true, 1, $false, "diferentes grupos.^lAvaliação das competições", 2) | Out-Null

# --- English "Programa" paragraph (detailed) ---
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Fields of action.Physics as a conceptual science", $true, $false, $false, $false, $false, $true, 1, $false, "Fields of action.^lPhysics as a conceptual science", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("experiments in Physics.Basic engineering concepts", $true, $false, $false, $false, $false, $true, 1, $false, "experiments in Physics.^lBasic engineering concepts", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("competences of an engineer.Development of a thematic project", $true, $false, $false, $false, $false, $true, 1, $false, "competences of an engineer.^lDevelopment of a thematic project", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Physical Engineering.Competition between projects", $true, $false, $false, $false, $false, $true, 1, $false, "Physical Engineering.^lCompetition between projects", 2) | Out-Null

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("from different groups.Evaluation of competitions", $true, $false, $false, $false, $false, $true, 1, $false, "from different groups.^lEvaluation of competitions", 2) | Out-Null
